# Generate Report for Handoff
# Swap the two source files' row order (ce6dfd37... now first, 5e99a1ad... second)
# and update statuses / handoff timestamps to reflect the new handoff report.

$wb = $excel.ActiveWorkbook

$ceFile = "ce6dfd37-98aa-48f9-a1f0-4da7f212bd8e.md"
$fiveFile = "5e99a1ad-3372-4ef2-a627-e0af3e2822bc.md"
$ceBase = "ce6dfd37-98aa-48f9-a1f0-4da7f212bd8e"
$fiveBase = "5e99a1ad-3372-4ef2-a627-e0af3e2822bc"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/1a0b6ea769d2ae440159cde039820858d321f075/e2e/"
$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/329594adc98be5e460d8949e5beea31f0d183974/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63c1daa2754f1cc35ab36e12d6b2654c101d9123/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $ceFile
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-13-11 16:13:00"

$wsOverview.Range("A3").Value = $fiveFile
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-13-11 16:13:00"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBase + $ceFile), "", "", $ceFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($mdBase + $fiveFile), "", "", $fiveFile)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhCeXlf = $ceBase + ".3f5327dbef68f717574d8bab9fff0402c54f0124.zh-cn.xlf"
$zhFiveXlf = $fiveBase + ".ea5881708952fc4c85cae9fd237d8870bf8bb766.zh-cn.xlf"

$wsZh.Range("A2").Value = $ceFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("D2").Value = $zhCeXlf
$wsZh.Range("E2").Value = "2016-03-11 16:12:57"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $fiveFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $zhFiveXlf
$wsZh.Range("E3").Value = "2016-03-11 16:13:00"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdBase + $ceFile), "", "", $ceFile)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), ($mdBase + $ceFile), "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), ($zhBase + $zhCeXlf), "", "", $zhCeXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($mdBase + $fiveFile), "", "", $fiveFile)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), ($mdBase + $fiveFile), "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), ($zhBase + $zhFiveXlf), "", "", $zhFiveXlf)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deCeXlf = $ceBase + ".3f5327dbef68f717574d8bab9fff0402c54f0124.de-de.xlf"
$deFiveXlf = $fiveBase + ".ea5881708952fc4c85cae9fd237d8870bf8bb766.de-de.xlf"

$wsDe.Range("A2").Value = $ceFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("D2").Value = $deCeXlf
$wsDe.Range("E2").Value = "2016-03-11 16:12:28"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $fiveFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $deFiveXlf
$wsDe.Range("E3").Value = "2016-03-11 16:13:00"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdBase + $ceFile), "", "", $ceFile)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), ($mdBase + $ceFile), "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), ($deBase + $deCeXlf), "", "", $deCeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($mdBase + $fiveFile), "", "", $fiveFile)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), ($mdBase + $fiveFile), "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), ($deBase + $deFiveXlf), "", "", $deFiveXlf)
